$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.146.17"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "2.252.92"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "98.82"
$ws.Range("E5").Value = "  +19.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "272.28"
$ws.Range("E6").Value = "  +5.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.648"
$ws.Range("E9").Value = "  +9.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.41"
$ws.Range("E10").Value = "  +8.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0955"
$ws.Range("E11").Value = "  +4.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.45"
$ws.Range("E12").Value = "  +19.34%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("E14").Value = "  +7.64%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.831"
$ws.Range("E15").Value = "  +6.77%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.247.51"
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "2.283.19"
$ws.Range("E17").Value = "  -9.39%  "
$ws.Range("D18").Value = "44.113.94"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000105"
$ws.Range("E19").Value = "  +3.71%  "
$ws.Range("E20").Value = "  +6.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.45"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.27"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.34"
$ws.Range("E24").Value = "  +6.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.69"
$ws.Range("E25").Value = "  +10.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +13.65%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.76"
$ws.Range("E28").Value = "  +7.97%  "
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.51"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.59"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0925"
$ws.Range("E32").Value = "  +7.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.09"
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.64"
$ws.Range("E34").Value = "  +6.34%  "
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0355"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.79"
$ws.Range("E39").Value = "  +34.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.244"
$ws.Range("E40").Value = "  +23.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.13"
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.18"
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.46"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.48"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("E45").Value = "  +3.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.55"
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("E47").Value = "  +3.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.17"
$ws.Range("E48").Value = "  +5.63%  "
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.444"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").Value = "2.475.92"
$ws.Range("E51").Value = "  +2.76%  "
